$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Nueva tabla "ProductoDeuda" (filas 22-28) ---
$ws.Range("F22").Value = "ProductoDeuda"

$ws.Range("F23").Value = "id_credito"
$ws.Range("G23").Value = "nombre"
$ws.Range("H23").Value = "cantidad"
$ws.Range("I23").Value = "precio"

$ws.Range("F24").Value = 1
$ws.Range("G24").Value = "coca"
$ws.Range("H24").Value = 1
$ws.Range("I24").Value = 18

$ws.Range("F25").Value = 1
$ws.Range("G25").Value = "pepsi "
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = 20

$ws.Range("F26").Value = 1
$ws.Range("G26").Value = "gansito"
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 30

$ws.Range("F27").Value = 2
$ws.Range("G27").Value = "desodorante"
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 55

$ws.Range("F28").Value = 3
$ws.Range("G28").Value = "coca"
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 35

# --- Nueva fila de credito (fila 19) ---
$ws.Range("A19").Value = 3
$ws.Range("C19").Value = "30_julio_2022"
$ws.Range("D19").Value = 1
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 150

# --- Se quita la columna id_cliente de la tabla de credito (ventana de abono) ---
$ws.Range("E16").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("E18").ClearContents()

# --- Selección final ---
$ws.Range("D23").Select()
